# "Adding Scenario 5 Test Case 3"
#
# This adds a new "FourthPage" worksheet (a J-column form identical in
# layout to "SecondPage" but starting at Fname/Lname like "FirstPage"),
# makes it the active/selected tab, and updates "FirstPage"'s selection
# to match the no-longer-active state it is left in.

$wb = $excel.ActiveWorkbook

# --- FirstPage: it is no longer the tabSelected sheet once FourthPage is
#     activated below; just update its lingering selection rectangle. ---
$firstPage = $wb.Worksheets.Item("FirstPage")
$firstPage.Range("D2:K2").Select()

# --- Build FourthPage. SecondPage already has the exact target column
#     layout (A:J, same widths/styles/hyperlink), so copy it and tweak
#     the two header/data cells that differ (Fname/Lname vs Lname/JFunc). ---
$template = $wb.Worksheets.Item("SecondPage")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)

$fourthPage = $wb.Worksheets.Item($wb.Worksheets.Count)
$fourthPage.Name = "FourthPage"

$fourthPage.Range("A1").Value = "Fname"
$fourthPage.Range("B1").Value = "Lname"
$fourthPage.Range("A2").Value = "Selenium"
$fourthPage.Range("B2").Value = "Automaition"

# --- Leave FourthPage as the selected/active tab with its own selection. ---
$fourthPage.Range("I13").Select()
$fourthPage.Activate()
